$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Jugadores")
$ws2 = $wb.Worksheets.Item("Entrenadores")

# --- Fix "Jugadores": swap rows 3 and 4 (Lorenzo Albarracin <-> Agustin Lagos) ---
# Use a scratch row far below the used range as temporary holding space.
$scratch = 1000

$ws1.Rows.Item(3).Copy()
$ws1.Rows.Item($scratch).PasteSpecial(-4104)   # xlPasteAll

$ws1.Rows.Item(4).Copy()
$ws1.Rows.Item(3).PasteSpecial(-4104)

$ws1.Rows.Item($scratch).Copy()
$ws1.Rows.Item(4).PasteSpecial(-4104)

$ws1.Rows.Item($scratch).Delete()
$excel.CutCopyMode = 0

# Row 3 now holds Agustin Lagos' data, row 4 holds Lorenzo Albarracin's data.

# Correct the shield / player photo filenames and move the misplaced
# "Aspectos_Tecnicos" text (was typed one column over, in AZ) into BB.
$ws1.Range("AX3").Value = "Tallers"
$ws1.Range("AY3").Value = "Facundo_Bernal"
$ws1.Range("BB3").Value = $ws1.Range("AZ3").Value2
$ws1.Range("AZ3").ClearContents()

$ws1.Range("AX4").Value = "Argentinos"
$ws1.Range("AY4").Value = "Diego_Armando_Maradona"
$ws1.Range("AZ4").ClearContents()

# Agustin Lagos' row needed a bit more room once the text landed in the right column.
$ws1.Rows.Item(3).RowHeight = 152

# --- Switch the active sheet/selection from "Entrenadores" to "Jugadores" ---
$ws1.Activate()
$ws1.Range("A4").Select()
